$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.363.97'
$ws.Cells.Item(2, 5).Value = '  -1.26%  '

$ws.Cells.Item(3, 4).Value = '2.519.21'
$ws.Cells.Item(3, 5).Value = '  -1.98%  '

$ws.Cells.Item(4, 5).Value = '  -0.02%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '303.03'

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '95.90'
$ws.Cells.Item(6, 5).Value = '  -1.46%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.583'
$ws.Cells.Item(7, 5).Value = '  +1.23%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.533'
$ws.Cells.Item(9, 5).Value = '  -2.79%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '36.36'
$ws.Cells.Item(10, 5).Value = '  +0.06%  '

$ws.Cells.Item(11, 5).Value = '  -0.29%  '

$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.112'
$ws.Cells.Item(12, 5).Value = '  -1.57%  '

$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '7.45'
$ws.Cells.Item(13, 5).Value = '  -2.62%  '

$ws.Cells.Item(14, 4).Value = '2.905.62'
$ws.Cells.Item(14, 5).Value = '  -2.00%  '

$ws.Cells.Item(15, 4).Value = '2.520.32'
$ws.Cells.Item(15, 5).Value = '  -1.95%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '14.99'
$ws.Cells.Item(16, 5).Value = '  +4.28%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.860'
$ws.Cells.Item(17, 5).Value = '  -2.96%  '

$ws.Cells.Item(18, 4).Value = '42.422.46'
$ws.Cells.Item(18, 5).Value = '  -1.18%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '12.84'
$ws.Cells.Item(19, 5).Value = '  -0.19%  '

$ws.Cells.Item(20, 5).Value = '  -2.81%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.43'
$ws.Cells.Item(21, 5).Value = '  -3.06%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '70.93'
$ws.Cells.Item(22, 5).Value = '  -1.50%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '250.18'
$ws.Cells.Item(23, 5).Value = '  -1.64%  '

$ws.Cells.Item(24, 5).Value = '  -1.93%  '

$ws.Cells.Item(25, 5).Value = '  -5.19%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '26.98'
$ws.Cells.Item(26, 5).Value = '  -6.40%  '

$ws.Cells.Item(27, 5).Value = '  +0.23%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.32'
$ws.Cells.Item(28, 5).Value = '  +9.56%  '

$ws.Cells.Item(29, 5).Value = '  -0.04%  '

$ws.Cells.Item(30, 5).Value = '  +1.20%  '

$ws.Cells.Item(31, 5).Value = '  -1.64%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '154.79'
$ws.Cells.Item(32, 5).Value = '  -0.26%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.32'
$ws.Cells.Item(33, 5).Value = '  -2.57%  '

$ws.Cells.Item(34, 5).Value = '  -2.41%  '

$ws.Cells.Item(35, 2).Value = 'Celestia'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '18.60'
$ws.Cells.Item(35, 5).Value = '  +1.58%  '

$ws.Cells.Item(36, 2).Value = 'WEMIXToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.62'
$ws.Cells.Item(36, 5).Value = '  -5.20%  '

$ws.Cells.Item(37, 5).Value = '  -4.86%  '

$ws.Cells.Item(38, 5).Value = '  +1.08%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '24.19'
$ws.Cells.Item(39, 5).Value = '  +4.66%  '

$ws.Cells.Item(40, 5).Value = '  -0.78%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.37'
$ws.Cells.Item(41, 5).Value = '  -1.63%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.82'
$ws.Cells.Item(42, 5).Value = '  -1.49%  '

$ws.Cells.Item(43, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.00'
$ws.Cells.Item(43, 5).Value = '  +0.06%  '

$ws.Cells.Item(44, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.02'
$ws.Cells.Item(44, 5).Value = '  -2.74%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0298'
$ws.Cells.Item(45, 5).Value = '  -3.82%  '

$ws.Cells.Item(46, 4).Value = '2.028.76'
$ws.Cells.Item(46, 5).Value = '  -2.28%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '84.52'
$ws.Cells.Item(47, 5).Value = '  -0.84%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '8.93'
$ws.Cells.Item(48, 5).Value = '  -2.90%  '

$ws.Cells.Item(49, 4).Value = '2.765.18'
$ws.Cells.Item(49, 5).Value = '  -2.04%  '

$ws.Cells.Item(50, 2).Value = 'Aave'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '101.38'
$ws.Cells.Item(50, 5).Value = '  -4.93%  '

$ws.Cells.Item(51, 2).Value = 'Algorand'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.188'
$ws.Cells.Item(51, 5).Value = '  -1.14%  '
